$wb = $excel.ActiveWorkbook

# --- Insert a new "Texas Notes" worksheet between "About" and "PDiCECpDoC" ---
$about = $wb.Worksheets.Item("About")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $about)
$newSheet.Name = "Texas Notes"

# Re-fetch the PDiCECpDoC worksheet reference (object references can go stale
# after inserting a new sheet into the collection)
$pdi = $wb.Worksheets.Item("PDiCECpDoC")

# --- Populate the Texas Notes sheet ---
$newSheet.Range("A1").Value = "A more recent study from DNVGL "

$newSheet.Range("A2").Value = "https://www.dnvgl.com/feature/carbon-capture-storage-ccs.html"
$newSheet.Hyperlinks.Add($newSheet.Range("A2"), "https://www.dnvgl.com/feature/carbon-capture-storage-ccs.html") | Out-Null

$newSheet.Range("A3").Value = "assumes a learning rate of 15-20% - closer to what we see in the renewables and storage industries."

$newSheet.Range("A5").Value = "This suggests that the learning rate may be higher than the 2013 report used by EI. "
$newSheet.Range("A6").Value = "So, we can take an average of some of these values just to be conservative"

$newSheet.Range("B7").Value = "average learning rate"

$newSheet.Range("A8").Value = "2018 DNVGL"
$newSheet.Range("A8").HorizontalAlignment = -4131
$newSheet.Range("B8").Formula = "=AVERAGE(0.15, 0.2)"

$newSheet.Range("A9").Value = "2013 CRS report"
$newSheet.Range("A9").HorizontalAlignment = -4131
$newSheet.Range("B9").Value = 0.13

$newSheet.Range("A10").Value = "average"
$newSheet.Range("A10").HorizontalAlignment = -4152
$newSheet.Range("B10").Formula = "=AVERAGE(B8:B9)"
$newSheet.Range("B10").Interior.Color = 65535

# Column width for column A on the new sheet
$newSheet.Columns.Item(1).ColumnWidth = 16.83

# --- Update PDiCECpDoC!B2 to reference the new average ---
$pdi.Range("B2").Formula = "='Texas Notes'!B10"

# --- View / selection state ---
$about.Activate()
$about.Range("E24").Select() | Out-Null

$newSheet.Activate()
$newSheet.Range("G13").Select() | Out-Null

$pdi.Activate()
$pdi.Range("C7").Select() | Out-Null
